# Updated mineral properties with SA, calculated SA normalized Kd
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compilation")

# New headers for columns G and H
$ws.Range("H2").Value = "K normalized by SA"
$ws.Range("G2").Value = "SA"

# Surface area (SA) values per mineral/electrolyte group
# Ferrihydrite, 10 mM NaCl (rows 3-6)
$ws.Range("G3").Value = 382.9
$ws.Range("G4").Value = 382.9
$ws.Range("G5").Value = 382.9
$ws.Range("G6").Value = 382.9

# Goethite, 10 mM NaCl (rows 7-10) - row 7 Kd is N/A, leave SA blank
$ws.Range("G7").HorizontalAlignment = 1
$ws.Range("G8").Value = 146.46
$ws.Range("G9").Value = 146.46
$ws.Range("G10").Value = 146.46

# Rows with no SA data but still part of the new column range (kept blank, styled)
$ws.Range("G11").HorizontalAlignment = 1
$ws.Range("G15").HorizontalAlignment = 1

# Sodium Montmorillonite, 10 mM NaCl Experimental (rows 20-23)
$ws.Range("G20").Value = 50.162
$ws.Range("G20").HorizontalAlignment = 1
$ws.Range("G21").Value = 50.162
$ws.Range("G21").HorizontalAlignment = 1
$ws.Range("G22").Value = 50.162
$ws.Range("G22").HorizontalAlignment = 1
$ws.Range("G23").Value = 50.162
$ws.Range("G23").HorizontalAlignment = 1

# K normalized by SA = Kd (column E) / SA (column G), for rows 4-23 (shared formula block)
$ws.Range("H4:H23").Formula = "=E4/G4"
$ws.Range("H7").ClearContents()
$ws.Range("H11:H19").ClearContents()

# H3 is a standalone (non-shared) formula
$ws.Range("H3").Formula = "=E3/G3"

$ws.Range("H5").Select()
